$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.078.72'
$ws.Range("D3").Value = '2.107.24'
$ws.Range("E3").Value = '  -0.66%  '
$ws.Range("E4").Value = '  -0.72%  '
$ws.Range("D5").Value = '''345.99'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("D7").Value = '''0.5183'
$ws.Range("E7").Value = '  -1.84%  '
$ws.Range("D8").Value = '''0.4437'
$ws.Range("E8").Value = '  -2.64%  '
$ws.Range("D9").Value = '''0.09411'
$ws.Range("E9").Value = '  +3.20%  '
$ws.Range("D10").Value = '''52.54'
$ws.Range("E10").Value = '  -2.90%  '
$ws.Range("D11").Value = '''1.179'
$ws.Range("E11").Value = '  +0.31%  '
$ws.Range("D12").Value = '''25.40'
$ws.Range("E12").Value = '  +3.66%  '
$ws.Range("D13").Value = '2.110.14'
$ws.Range("E13").Value = '  -0.56%  '
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("D15").Value = '''8.160'
$ws.Range("E15").Value = '  +0.78%  '
$ws.Range("D16").Value = '''99.94'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").Value = '''0.00001169'
$ws.Range("E17").Value = '  -0.18%  '
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").Value = '''20.78'
$ws.Range("E19").Value = '  +6.17%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = '''1.006'
$ws.Range("E21").Value = '  -0.74%  '
$ws.Range("D22").Value = '''6.236'
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("D23").Value = '30.174.34'
$ws.Range("E23").Value = '  -1.94%  '
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").Value = '''2.343'
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").Value = '2.358.97'
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("E27").Value = '  -1.81%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '''164.23'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").Value = '''2.554'
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").Value = '''133.98'
$ws.Range("E30").Value = '  -0.85%  '
$ws.Range("D31").Value = '''1.171'
$ws.Range("E31").Value = '  -3.26%  '
$ws.Range("D32").Value = '''0.1062'
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("D33").Value = '''1.643'
$ws.Range("E33").Value = '  +0.27%  '
$ws.Range("D34").Value = '''6.269'
$ws.Range("E34").Value = '  -2.35%  '
$ws.Range("D35").Value = '''3.953'
$ws.Range("E35").Value = '  -0.10%  '
$ws.Range("D36").Value = '''6.241'
$ws.Range("E36").Value = '  +4.18%  '
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("E38").Value = '  -3.72%  '
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("D40").Value = '''0.2294'
$ws.Range("E40").Value = '  -1.46%  '
$ws.Range("D41").Value = '''0.7008'
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").Value = '''12.58'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("E43").Value = '  +4.08%  '
$ws.Range("D44").Value = '''0.6724'
$ws.Range("E44").Value = '  +3.77%  '
$ws.Range("D45").Value = '''14.24'
$ws.Range("E45").Value = '  -6.26%  '
$ws.Range("D46").Value = '''2.294'
$ws.Range("E46").Value = '  -0.86%  '
$ws.Range("D47").Value = '''3.641'
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("D48").Value = '''0.00000000352'
$ws.Range("E48").Value = '  -4.94%  '
$ws.Range("D49").Value = '''1.224'
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("D50").Value = '''82.91'
$ws.Range("E50").Value = '  -0.34%  '
